# Re-run SGNN to annotate dialog acts following clean up work to the original transcripts.
# This updates the DAMSLTag (column I) and DialogAct (column J) values for a set of rows
# on the active worksheet, reflecting the re-annotated dialog act labels.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    13  = @("b", "Acknowledge (Backchannel)")
    14  = @("sd", "Statement-non-opinion")
    15  = @("aa", "Agree/Accept")
    19  = @("%", "Uninterpretable")
    22  = @("%", "Uninterpretable")
    23  = @("sd", "Statement-non-opinion")
    26  = @("sd", "Statement-non-opinion")
    27  = @("aa", "Agree/Accept")
    29  = @("aa", "Agree/Accept")
    31  = @("aa", "Agree/Accept")
    32  = @("aa", "Agree/Accept")
    33  = @("aa", "Agree/Accept")
    34  = @("aa", "Agree/Accept")
    42  = @("sv", "Statement-opinion")
    43  = @("sv", "Statement-opinion")
    74  = @("sv", "Statement-opinion")
    85  = @("sv", "Statement-opinion")
    89  = @("aa", "Agree/Accept")
    90  = @("sd", "Statement-non-opinion")
    96  = @("sd", "Statement-non-opinion")
    118 = @("aa", "Agree/Accept")
    124 = @("sd", "Statement-non-opinion")
    144 = @("sd", "Statement-non-opinion")
    156 = @("%", "Uninterpretable")
    164 = @("aa", "Agree/Accept")
    168 = @("%", "Uninterpretable")
    173 = @("ba", "Appreciation")
}

foreach ($rowNum in $updates.Keys) {
    $values = $updates[$rowNum]
    $ws.Range("I$rowNum").Value = $values[0]
    $ws.Range("J$rowNum").Value = $values[1]
}
